# Toevoegen nieuwe wijken najaar 2022
# Update the "aantal" (D) column from 2 to 1 for a set of rows, and
# rename a number of wijk-codes in column E to reflect new subdivided
# wijken (ggw7_type sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only the D-column (count) changes from 2.00 to 1.00
$dOnlyRows = @(93, 101, 102, 131, 140, 146, 148)

foreach ($r in $dOnlyRows) {
    $ws.Cells.Item($r, 4).Value = 1
}

# Rows where both the D-column changes (2.00 -> 1.00) and the E-column
# (wijkcode) text changes.
$eChanges = @{
    53  = "13013BER"
    112 = "23086B0"
    161 = "31033CENT"
    171 = "33016MEB"
    198 = "36007INGNO"
}

foreach ($r in $eChanges.Keys) {
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = $eChanges[$r]
}
